$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Vendor"
$ws.Range("J1").ClearContents()
